# Replace the data table on Sheet1 with the new, corrected dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$header = @("MetalShape", "MetalNumber", "Validity_Check_Value_of_Translation_Error", "Validity_Check_Value_of_Rotation_Error")
for ($c = 1; $c -le 4; $c++) {
  $ws.Cells.Item(1, $c).Value = $header[$c - 1]
}

$rows = @(
  @("sheet", "Control", 0.0000015459890823990422, 0.000028067641897829744),
  @("sheet", "LC Steel", 0.0000069627762128385077, 0.000060956958699467337),
  @("sheet", "304 SS", 0.0000038761120837528648, 0.000024410464528341806),
  @("sheet", "6061 Al", 0.000006737042247821481, 0.000022569920094691914),
  @("sheet", "Copper", 0.0000046053454519269281, 0.000040069836989257042),
  @("solid", "Ti Grade 5", 0.0000019637988763980749, 0.0000088177994831109279),
  @("solid", "Copper", 0.00000050093775271640976, 0.0000062043459560144227),
  @("solid", "Control", 0.0000011764180270163276, 0.0000098327321610156074),
  @("solid", "LC Steel", 0.0000020776531428105052, 0.000019707496771217918),
  @("solid", "416 SS", 0.0000035072702529884472, 0.000014327881731271386),
  @("solid", "304 SS", 0.0000031165359883408649, 0.000037015602558869957),
  @("solid", "6061 Al", 0.0000051375673733102401, 0.000029575898396514673),
  @("solid", "Ti Grade 5", 0.00000036446812295075986, 0.000014685192566197794),
  @("solid", "Copper", 0.0000021913802358796313, 0.000014689386105367674),
  @("solid", "Ti Grade 5", 0.0000002340384168052887, 0.000035065785209123983),
  @("solid", "Copper", 0.0000031594882232792449, 0.00003660722105154054)
)

$r = 2
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}
